$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 77, pushing the existing data (rows 77-99) down to rows 79-101.
$ws.Rows.Item(77).EntireRow.Insert()
$ws.Rows.Item(78).EntireRow.Insert()

# Fill the first new row (77) with the new "Primera" price entry for date serial 44468.
$ws.Cells.Item(77, 1).Value = 11
$ws.Cells.Item(77, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(77, 3).Value = "Bíobío"
$ws.Cells.Item(77, 4).Value = 44468
$ws.Cells.Item(77, 5).Value = 8
$ws.Cells.Item(77, 6).Value = "Fruta"
$ws.Cells.Item(77, 7).Value = 100101
$ws.Cells.Item(77, 8).Value = "Berries"
$ws.Cells.Item(77, 9).Value = 100101007
$ws.Cells.Item(77, 10).Value = "Kiwi"
$ws.Cells.Item(77, 11).Value = "Hayward"
$ws.Cells.Item(77, 12).Value = "Primera"
$ws.Cells.Item(77, 13).Value = 200
$ws.Cells.Item(77, 14).Value = 12000
$ws.Cells.Item(77, 15).Value = 13000
$ws.Cells.Item(77, 16).Value = 12500
$ws.Cells.Item(77, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(77, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(77, 19).Value = 694
$ws.Cells.Item(77, 20).Value = 18

# Fill the second new row (78) with the new "Segunda" price entry for the same date serial 44468.
$ws.Cells.Item(78, 1).Value = 11
$ws.Cells.Item(78, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(78, 3).Value = "Bíobío"
$ws.Cells.Item(78, 4).Value = 44468
$ws.Cells.Item(78, 5).Value = 8
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100101
$ws.Cells.Item(78, 8).Value = "Berries"
$ws.Cells.Item(78, 9).Value = 100101007
$ws.Cells.Item(78, 10).Value = "Kiwi"
$ws.Cells.Item(78, 11).Value = "Hayward"
$ws.Cells.Item(78, 12).Value = "Segunda"
$ws.Cells.Item(78, 13).Value = 100
$ws.Cells.Item(78, 14).Value = 10000
$ws.Cells.Item(78, 15).Value = 10000
$ws.Cells.Item(78, 16).Value = 10000
$ws.Cells.Item(78, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(78, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 19).Value = 556
$ws.Cells.Item(78, 20).Value = 18
